# Update census data
# The underlying survey changed (IDP disability register -> unified targeted
# social-assistance disability register), so:
#   - the report title (row 1) gets new wording
#   - the old single "Number of disability persons" row is replaced by two
#     rows pulled from the new source workbook: "family with disabilities
#     Persons" and "disabilities Persons", each with refreshed figures
#   - the source note (now row 6) and the year header keep their wording.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: insert a new row ABOVE the old data row (old row 4) first, while
# every other row is still in its original place, so we can safely copy
# formats from stable, untouched cells.
#   after insert:
#     row1 = title                (unchanged position)
#     row2 = "(End of year..)"    (unchanged position)
#     row3 = year header          (unchanged position)
#     row4 = brand new blank row
#     row5 = old data row (shifted down 1) - still has the ORIGINAL format
#     row6 = old source row (shifted down 1) - still has the ORIGINAL format
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).Insert()

# Grab the still-original row 5 (old data row) formatting for the new row 4
# label + data cells before we overwrite row 5's own content. Use the plain
# "B" data-cell format (no border, no explicit alignment) uniformly for all
# of the new numeric cells so they all share one style, matching how the
# pasted-in replacement rows come in from the source workbook.
$ws.Range("A5").Copy()
$ws.Range("A4").PasteSpecial(-4122)

$ws.Range("B5").Copy()
$ws.Range("B4:I4").PasteSpecial(-4122)
$ws.Range("B5:H5").PasteSpecial(-4122)
$ws.Range("I5").PasteSpecial(-4122)

$ws.Rows.Item(4).RowHeight = 24.75

# ---------------------------------------------------------------------------
# Step 2: title row (row 1)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Tsalenjikha Municipality"
$ws.Range("A1:I1").Merge()
$ws.Range("A1:I1").HorizontalAlignment = -4108
$ws.Range("A1:I1").VerticalAlignment = -4108
$ws.Range("A1:I1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 51

# ---------------------------------------------------------------------------
# Step 3: row 4 "family with disabilities Persons" (new row, formats already
# copied from the original data row above)
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "family with disabilities Persons "
$ws.Range("A4").Borders.Item(9).LineStyle = -4142

$ws.Range("B4").Value = 597
$ws.Range("C4").Value = 572
$ws.Range("D4").Value = 567
$ws.Range("E4").Value = 587
$ws.Range("F4").Value = 576
$ws.Range("G4").Value = 563
$ws.Range("H4").Value = 562
$ws.Range("I4").Value = 571

# ---------------------------------------------------------------------------
# Step 4: row 5 "disabilities Persons" - reuses the untouched "(End of
# year..)" label look (theme text colour, Arial 10, white fill) with a
# bottom rule instead of the old top+bottom pair, plus the existing numeric
# style for the figures.
# ---------------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = "disabilities Persons "
$ws.Range("A5").HorizontalAlignment = -4131
$ws.Range("A5").VerticalAlignment = -4108
$ws.Range("A5").WrapText = $true
$ws.Range("A5").Borders.Item(9).LineStyle = 1
$ws.Rows.Item(5).RowHeight = 21

$ws.Range("B5").Value = 671
$ws.Range("C5").Value = 649
$ws.Range("D5").Value = 639
$ws.Range("E5").Value = 655
$ws.Range("F5").Value = 639
$ws.Range("G5").Value = 622
$ws.Range("H5").Value = 617
$ws.Range("I5").Value = 632
$ws.Range("I5").Borders.Item(9).LineStyle = 1

# ---------------------------------------------------------------------------
# Step 5: the blank A3 cell above the year header picked up the pasted-in
# workbook's default body font (Sylfaen) when the replacement rows were
# brought in.
# ---------------------------------------------------------------------------
$ws.Range("A3").Font.Name = "Sylfaen"
$ws.Range("A3").Font.Size = 11

"done"
